# Adds the "Hacer Hadoken" task (row 6) and its "Modulo player 2" sound task
# (row 7) to the task-time tracking sheet, as per commit:
# "Añadido Hadoken y su respectivo sonido."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: Hacer Hadoken / Jose / 2 h / 2h
$ws.Range("A6").Value = "Hacer Hadoken"
$ws.Range("B6").Value = "Jose"
$ws.Range("C6").Value = "2 h"
$ws.Range("D6").Value = "2h"

# Row 7: Modulo player 2 / Jose, Ferran / 1 h / 1h
$ws.Range("A7").Value = "Modulo player 2"
$ws.Range("B7").Value = "Jose, Ferran"
$ws.Range("C7").Value = "1 h"
$ws.Range("D7").Value = "1h"

# Row 7 now matches row 6's (auto) height of 13.8 points.
$ws.Rows.Item(7).RowHeight = 13.8

# Move the active selection to C8, matching the post-edit cursor position.
[void]$ws.Range("C8").Select()
